$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Biblioteca" (library) values in column D for rows whose
# double-degree library text was previously truncated / mismatched.
# (Order matters for how new shared-string entries get appended.)
$ws.Range("D6").Value  = "B. Derecho + B. Políticas y Sociolog."
$ws.Range("D16").Value = "B. Ciencias + B. Informática y Telecom."
$ws.Range("D15").Value = " B. Informática y Telecom. + B. Económicas y Empres."
$ws.Range("D9").Value  = "B. Filosofía y Letras A + B. Eduación"
$ws.Range("D10").Value = "B. Filosofía y Letras A + B. Eduación"
$ws.Range("D7").Value  = "B. Económicas y Empres. + B. Politécnica"
$ws.Range("D13").Value = "B. Económicas y Empres. + B. Politécnica"
$ws.Range("D14").Value = "B. Económicas y Empres. + B. Politécnica"
$ws.Range("D4").Value  = "B. Políticas y Sociolog. + B. Colegio Máximo"

# Row 40 (Grado en Edificación) previously had no library assigned.
$ws.Range("D40").Value = "B. Politécnica"

# Restore the selection recorded in the sheet view.
$ws.Range("A48").Select()
